# Generate Report for Archive
# Updates the localization status report: two files ("0b2b01d3-...md" and
# "1ad09cfa-...md") have progressed from "Ready for handoff" back to
# "In Translation". Update the Status column on every sheet that tracks it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5").Value = "In Translation"
$overview.Range("F5").Value = "In Translation"
$overview.Range("E6").Value = "In Translation"
$overview.Range("F6").Value = "In Translation"

# --- zh-cn sheet: column C holds the Status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5").Value = "In Translation"
$zhcn.Range("C6").Value = "In Translation"

# --- de-de sheet: column C holds the Status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5").Value = "In Translation"
$dede.Range("C6").Value = "In Translation"
